$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "last row" (row 77) had a date-only number format to mark it
# as the latest entry. Since we're appending a new row, row 77 reverts to
# the regular datetime number format used by every other data row.
$ws.Cells.Item(77, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new daily row (day 45665) with its values.
$ws.Cells.Item(78, 1).Value = 45665
$ws.Cells.Item(78, 2).Value = 182
$ws.Cells.Item(78, 3).Value = 181
$ws.Cells.Item(78, 4).Value = 183

# The new last row gets the date-only number format that previously marked
# row 77 as the latest entry.
$ws.Cells.Item(78, 1).NumberFormat = "YYYY-MM-DD"
